$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 9999
$ws.Range("J10").Value = 9999
$ws.Range("L10").Value = 9999
$ws.Range("N10").Value = -10585
$ws.Range("H18").Value = 9998.625
$ws.Range("J18").Value = 818
$ws.Range("L18").Value = 818
$ws.Range("N18").Value = -1386
$ws.Range("H40").Value = 4462.5
$ws.Range("I40").Value = 1600
$ws.Range("J40").Value = 5416.6665
$ws.Range("K40").Value = 1600
$ws.Range("L40").Value = 5416.6665
$ws.Range("M40").Value = -1425
$ws.Range("N40").Value = -5766.6665
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H64").Value = 11999.5
$ws.Range("I64").Value = 3999
$ws.Range("J64").Value = 20000
$ws.Range("K64").Value = 3999
$ws.Range("L64").Value = 20000
$ws.Range("M64").Value = -3751
$ws.Range("N64").Value = -20496
$ws.Range("H67").Value = 11999.5
$ws.Range("I67").Value = 3999
$ws.Range("J67").Value = 20000
$ws.Range("K67").Value = 3999
$ws.Range("L67").Value = 20000
$ws.Range("M67").Value = -3141
$ws.Range("N67").Value = -21716
$ws.Range("H70").Value = 8481.6
$ws.Range("I70").Value = 1479.1666
$ws.Range("J70").Value = 11482.643
$ws.Range("K70").Value = 4437.4998
$ws.Range("L70").Value = 34447.929
$ws.Range("M70").Value = -4167.4998
$ws.Range("N70").Value = -34987.929
$ws.Range("H73").Value = 8481.6
$ws.Range("I73").Value = 1479.1666
$ws.Range("J73").Value = 11482.643
$ws.Range("K73").Value = 4437.4998
$ws.Range("L73").Value = 34447.929
$ws.Range("M73").Value = -3501.4998
$ws.Range("N73").Value = -36319.929
$ws.Range("H82").Value = 1637.4
$ws.Range("I82").Value = 1637.4
$ws.Range("K82").Value = 4912.200000000001
$ws.Range("M82").Value = -4506.200000000001
$ws.Range("H85").Value = 1637.4
$ws.Range("I85").Value = 1637.4
$ws.Range("K85").Value = 4912.200000000001
$ws.Range("M85").Value = -3508.200000000001
$ws.Range("H116").Value = 2222.2222
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 2000
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 1442
$ws.Range("N116").Value = -9884

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 2750
$ws.Range("H61").Value = 3500
$ws.Range("I61").Value = 3500
$ws.Range("K61").Value = 3500
$ws.Range("M61").Value = -3288
$ws.Range("H88").Value = 5248.5
$ws.Range("I88").Value = 2999.5
$ws.Range("J88").Value = 7497.5
$ws.Range("K88").Value = 2999.5
$ws.Range("L88").Value = 7497.5
$ws.Range("M88").Value = -2593.5
$ws.Range("N88").Value = -8309.5
$ws.Range("H91").Value = 5248.5
$ws.Range("I91").Value = 2999.5
$ws.Range("J91").Value = 7497.5
$ws.Range("K91").Value = 2999.5
$ws.Range("L91").Value = 7497.5
$ws.Range("M91").Value = -1595.5
$ws.Range("N91").Value = -10305.5
$ws.Range("H136").Value = 3500
$ws.Range("I136").Value = 3500
$ws.Range("K136").Value = 10500
$ws.Range("M136").Value = -7950

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 971.2857
$ws.Range("I86").Value = 859.8
$ws.Range("J86").Value = 1250
$ws.Range("K86").Value = 859.8
$ws.Range("L86").Value = 1250
$ws.Range("M86").Value = 263.2
$ws.Range("N86").Value = -3496
$ws.Range("H89").Value = 971.2857
$ws.Range("I89").Value = 859.8
$ws.Range("J89").Value = 1250
$ws.Range("K89").Value = 4299
$ws.Range("L89").Value = 6250
$ws.Range("M89").Value = 1317
$ws.Range("N89").Value = -17482

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 6800
$ws.Range("J3").Value = 6800
$ws.Range("L3").Value = 6800
$ws.Range("N3").Value = -7026
$ws.Range("H22").Value = 106.666664
$ws.Range("I22").Value = 106.666664
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 106.666664
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 243.333336
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 5299.375
$ws.Range("I31").Value = 1390.8
$ws.Range("J31").Value = 7076
$ws.Range("K31").Value = 1390.8
$ws.Range("L31").Value = 7076
$ws.Range("M31").Value = -1095.8
$ws.Range("N31").Value = -7666
$ws.Range("H34").Value = 5299.375
$ws.Range("I34").Value = 1390.8
$ws.Range("J34").Value = 7076
$ws.Range("K34").Value = 1390.8
$ws.Range("L34").Value = 7076
$ws.Range("M34").Value = -1188.8
$ws.Range("N34").Value = -7480
$ws.Range("H50").Value = 18666.666
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 18666.666
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 18666.666
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -19916.666
$ws.Range("H51").Value = 17333.334
$ws.Range("J51").Value = 17333.334
$ws.Range("L51").Value = 17333.334
$ws.Range("N51").Value = -18805.334
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H60").Value = 16199.9
$ws.Range("J60").Value = 18714.285
$ws.Range("L60").Value = 18714.285
$ws.Range("N60").Value = -19736.285
$ws.Range("H61").Value = 17333.334
$ws.Range("J61").Value = 17333.334
$ws.Range("L61").Value = 17333.334
$ws.Range("N61").Value = -18029.334
$ws.Range("H110").Value = 99995
$ws.Range("J110").Value = 99995
$ws.Range("L110").Value = 99995
$ws.Range("N110").Value = -108175

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 462104.84
$ws.Range("J4").Value = 500499.5
$ws.Range("L4").Value = 1501498.5
$ws.Range("N4").Value = -1501722.5
$ws.Range("H11").Value = 50000000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H63").Value = 9000
$ws.Range("J63").Value = 9000
$ws.Range("L63").Value = 27000
$ws.Range("N63").Value = -28498
$ws.Range("H66").Value = 9000
$ws.Range("J66").Value = 9000
$ws.Range("L66").Value = 81000
$ws.Range("N66").Value = -88488
$ws.Range("H107").Value = 619.8
$ws.Range("I107").Value = 499.5
$ws.Range("K107").Value = 1498.5
$ws.Range("M107").Value = 421.5
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1408.3334
$ws.Range("I22").Value = 1400
$ws.Range("J22").Value = 1450
$ws.Range("K22").Value = 1400
$ws.Range("L22").Value = 1450
$ws.Range("M22").Value = -1105
$ws.Range("N22").Value = -2040
$ws.Range("H27").Value = 1408.3334
$ws.Range("I27").Value = 1400
$ws.Range("J27").Value = 1450
$ws.Range("K27").Value = 1400
$ws.Range("L27").Value = 1450
$ws.Range("M27").Value = -1293
$ws.Range("N27").Value = -1664
$ws.Range("H100").Value = 2851.3
$ws.Range("I100").Value = 2851.3
$ws.Range("K100").Value = 2851.3
$ws.Range("M100").Value = -2310.3
$ws.Range("H125").Value = 39000
$ws.Range("J125").Value = 39000
$ws.Range("L125").Value = 39000
$ws.Range("N125").Value = -48840
$ws.Range("H132").Value = 6279.1665
$ws.Range("I132").Value = 4535
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 13605
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -11075
$ws.Range("N132").Value = -50060

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 41306
$ws.Range("J4").Value = 2132
$ws.Range("L4").Value = 2132
$ws.Range("N4").Value = -2358
$ws.Range("H16").Value = 99995
$ws.Range("J16").Value = 99995
$ws.Range("L16").Value = 99995
$ws.Range("N16").Value = -100579
$ws.Range("H20").Value = 21870.334
$ws.Range("J20").Value = 21870.334
$ws.Range("L20").Value = 21870.334
$ws.Range("N20").Value = -22350.334
$ws.Range("H122").Value = 2914.1
$ws.Range("I122").Value = 2904
$ws.Range("J122").Value = 3005
$ws.Range("K122").Value = 8712
$ws.Range("L122").Value = 9015
$ws.Range("M122").Value = -6262
$ws.Range("N122").Value = -13915
$ws.Range("H132").Value = 4374.25
$ws.Range("I132").Value = 4165.6665
$ws.Range("K132").Value = 12496.9995
$ws.Range("M132").Value = -9966.999500000002
